# Auto-generated Excel COM-interop script applying the Tonberry_Profits diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 20
$ws.Range("H20").Value = 25000
$ws.Range("I20").Value = 25000
$ws.Range("K20").Value = 25000
$ws.Range("M20").Value = -24770

# Row 35
$ws.Range("H35").Value = 25000
$ws.Range("I35").Value = 25000
$ws.Range("K35").Value = 25000
$ws.Range("M35").Value = -24621

# Row 113
$ws.Range("H113").Value = 8888.4375
$ws.Range("I113").Value = 9282.333000000001
$ws.Range("J113").Value = 2980
$ws.Range("K113").Value = 9282.333000000001
$ws.Range("L113").Value = 2980
$ws.Range("M113").Value = -6028.333000000001
$ws.Range("N113").Value = -9488

# Row 133
$ws.Range("H133").Value = 64571.43
$ws.Range("J133").Value = 64571.43
$ws.Range("L133").Value = 64571.43
$ws.Range("N133").Value = -74691.42999999999

# Row 138
$ws.Range("H138").Value = 3611.6863
$ws.Range("J138").Value = 2939.2683
$ws.Range("L138").Value = 8817.804900000001
$ws.Range("N138").Value = -19097.8049

$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 20000074
$ws.Range("I5").Value = 25
$ws.Range("K5").Value = 25
$ws.Range("M5").Value = 87

# Row 32
$ws.Range("H32").Value = 10421.7705
$ws.Range("I32").Value = 8032.024
$ws.Range("J32").Value = 25679.385
$ws.Range("K32").Value = 8032.024
$ws.Range("L32").Value = 25679.385
$ws.Range("M32").Value = -7745.024
$ws.Range("N32").Value = -26253.385

# Row 45
$ws.Range("H45").Value = 4501167.5
$ws.Range("I45").Value = 8182448.5
$ws.Range("K45").Value = 8182448.5
$ws.Range("M45").Value = -8182071.5

# Row 74
$ws.Range("H74").Value = 850.54285
$ws.Range("I74").Value = 647.129
$ws.Range("J74").Value = 2427
$ws.Range("K74").Value = 647.129
$ws.Range("L74").Value = 2427
$ws.Range("M74").Value = 226.871
$ws.Range("N74").Value = -4175

# Row 77
$ws.Range("H77").Value = 850.54285
$ws.Range("I77").Value = 647.129
$ws.Range("J77").Value = 2427
$ws.Range("K77").Value = 3235.645
$ws.Range("L77").Value = 12135
$ws.Range("M77").Value = 1132.355
$ws.Range("N77").Value = -20871

# Row 132
$ws.Range("H132").Value = 1869.7391
$ws.Range("I132").Value = 1481.359
$ws.Range("J132").Value = 2374.6333
$ws.Range("K132").Value = 4444.076999999999
$ws.Range("L132").Value = 7123.8999
$ws.Range("M132").Value = -1914.076999999999
$ws.Range("N132").Value = -12183.8999

# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = $null

# Row 139
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").Value = $null

# Row 140
$ws.Range("H140").Value = 59100
$ws.Range("J140").Value = 59100
$ws.Range("L140").Value = 59100
$ws.Range("N140").Value = -69460

$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 20000074
$ws.Range("I4").Value = 25
$ws.Range("K4").Value = 25
$ws.Range("M4").Value = 90

# Row 107
$ws.Range("H107").Value = 2515.1667
$ws.Range("I107").Value = 2515.1667
$ws.Range("K107").Value = 2515.1667
$ws.Range("M107").Value = -595.1667000000002

# Row 134
$ws.Range("H134").Value = 14522.2
$ws.Range("I134").Value = 23362.4
$ws.Range("J134").Value = 5682
$ws.Range("K134").Value = 70087.20000000001
$ws.Range("L134").Value = 17046
$ws.Range("M134").Value = -67552.20000000001
$ws.Range("N134").Value = -22116

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 1208.8889
$ws.Range("I22").Value = 692.5
$ws.Range("K22").Value = 692.5
$ws.Range("M22").Value = -342.5

# Row 58
$ws.Range("H58").Value = 967955
$ws.Range("I58").Value = 5437187
$ws.Range("J58").Value = 1634.5405
$ws.Range("K58").Value = 5437187
$ws.Range("L58").Value = 1634.5405
$ws.Range("M58").Value = -5436984
$ws.Range("N58").Value = -2040.5405

# Row 134
$ws.Range("H134").Value = 1367.1818
$ws.Range("I134").Value = 1241.6
$ws.Range("J134").Value = 1759.625
$ws.Range("K134").Value = 3724.8
$ws.Range("L134").Value = 5278.875
$ws.Range("M134").Value = -1189.8
$ws.Range("N134").Value = -10348.875

# Row 136
$ws.Range("H136").Value = 967955
$ws.Range("I136").Value = 5437187
$ws.Range("J136").Value = 1634.5405
$ws.Range("K136").Value = 16311561
$ws.Range("L136").Value = 4903.6215
$ws.Range("M136").Value = -16309011
$ws.Range("N136").Value = -10003.6215

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 24017.334
$ws.Range("J131").Value = 29173.75
$ws.Range("L131").Value = 87521.25
$ws.Range("N131").Value = -97601.25

# Row 141
$ws.Range("H141").Value = 3734.4285
$ws.Range("I141").Value = 3734.4285
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 11203.2855
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -6023.2855
$ws.Range("N141").Value = $null

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 3306.8572
$ws.Range("I97").Value = 3306.8572
$ws.Range("K97").Value = 3306.8572
$ws.Range("M97").Value = -2810.8572

# Row 113
$ws.Range("H113").Value = 1866.3334
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1866.3334
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 1866.3334
$ws.Range("M113").Value = $null
$ws.Range("N113").Value = -6206.3334

# Row 122
$ws.Range("H122").Value = 1413.4615
$ws.Range("J122").Value = 1812.3334
$ws.Range("L122").Value = 5437.0002
$ws.Range("N122").Value = -10337.0002

# Row 132
$ws.Range("H132").Value = 942333.9399999999
$ws.Range("I132").Value = 1378597.1
$ws.Range("J132").Value = 2690.077
$ws.Range("K132").Value = 4135791.3
$ws.Range("L132").Value = 8070.231000000001
$ws.Range("M132").Value = -4133261.3
$ws.Range("N132").Value = -13130.231

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1202.6666
$ws.Range("I22").Value = 926.5714
$ws.Range("J22").Value = 1444.25
$ws.Range("K22").Value = 926.5714
$ws.Range("L22").Value = 1444.25
$ws.Range("M22").Value = -631.5714
$ws.Range("N22").Value = -2034.25

# Row 26
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").Value = $null

# Row 27
$ws.Range("H27").Value = 1202.6666
$ws.Range("I27").Value = 926.5714
$ws.Range("J27").Value = 1444.25
$ws.Range("K27").Value = 926.5714
$ws.Range("L27").Value = 1444.25
$ws.Range("M27").Value = -819.5714
$ws.Range("N27").Value = -1658.25

# Row 46
$ws.Range("H46").Value = 2051.5386
$ws.Range("J46").Value = 2565
$ws.Range("L46").Value = 2565
$ws.Range("N46").Value = -2941

# Row 136
$ws.Range("H136").Value = 1337.9286
$ws.Range("I136").Value = 968.6667
$ws.Range("J136").Value = 2002.6
$ws.Range("K136").Value = 2906.0001
$ws.Range("L136").Value = 6007.799999999999
$ws.Range("M136").Value = -356.0001000000002
$ws.Range("N136").Value = -11107.8

$ws = $wb.Worksheets.Item("WVR")
# Row 46
$ws.Range("H46").Value = 78357.25
$ws.Range("J46").Value = 78357.25
$ws.Range("L46").Value = 78357.25
$ws.Range("N46").Value = -78819.25

# Row 122
$ws.Range("H122").Value = 32083.814
$ws.Range("J122").Value = 4337.375
$ws.Range("L122").Value = 13012.125
$ws.Range("N122").Value = -17912.125

# Row 132
$ws.Range("H132").Value = 1097.1296
$ws.Range("I132").Value = 894.0454999999999
$ws.Range("J132").Value = 1990.7
$ws.Range("K132").Value = 2682.1365
$ws.Range("L132").Value = 5972.1
$ws.Range("M132").Value = -152.1364999999996
$ws.Range("N132").Value = -11032.1

# Row 134
$ws.Range("H134").Value = 78357.25
$ws.Range("J134").Value = 78357.25
$ws.Range("L134").Value = 235071.75
$ws.Range("N134").Value = -240141.75
